$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Insert a new row before row 13 (pushes current rows 13-16 down to 14-17)
$ws.Range("A13:G13").Insert(-4121)

# Copy format from row 12 (same B/C/D/E/G styles we need for the new row) into row 13
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$ws.Range("F13").Clear()

# Fill in the new row's content (new RepayEndDate field)
$ws.Cells.Item(13,1).Value2 = 5
$ws.Cells.Item(13,2).Value2 = "RepayEndDate"
$ws.Cells.Item(13,3).Value2 = "應繳截止日"
$ws.Cells.Item(13,4).Value2 = "Decimald"
$ws.Cells.Item(13,5).Value2 = 8
$ws.Cells.Item(13,7).Value2 = "2022-03-11新增`n原系統有此欄位`nTBYGYMP.YGEPDT"

# New note in G13 gets a highlighted (yellow) wrapped style
$g13 = $ws.Cells.Item(13,7)
$g13.Interior.Color = 65535
$g13.WrapText = $true

$ws.Rows(13).RowHeight = 59.4

# Renumber the SEQ column for the rows that got pushed down
$ws.Cells.Item(14,1).Value2 = 6
$ws.Cells.Item(15,1).Value2 = 7
$ws.Cells.Item(16,1).Value2 = 8
$ws.Cells.Item(17,1).Value2 = 9

# Update the sheet view: scrolled to row 7, selection on D12
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D12").Select()
